{"js": "// Replace the 25 \"NNN\u00f7N=\" division-problem prompts in the table with the\n// new set of prompts, matching each old value to its new value one-to-one.\n// Every old prompt text is unique in the document (confirmed against the\n// original OOXML), so an exact, case-sensitive whole-string search for each\n// old value finds exactly one run and can be swapped in place safely.\n//\n// NOTE on ordering: \"170\u00f79=\" is both an original value (last cell) and a\n// new value (first cell), so that pair is processed first -- before any\n// other insertion creates a second \"170\u00f79=\" in the document -- to keep\n// every search unambiguous (exactly one match at the time it runs).\nconst replacements = [\n  [\"170\u00f79=\", \"146\u00f75=\"],\n  [\"563\u00f74=\", \"170\u00f79=\"],\n  [\"174\u00f77=\", \"575\u00f72=\"],\n  [\"507\u00f72=\", \"399\u00f76=\"],\n  [\"733\u00f79=\", \"413\u00f79=\"],\n  [\"555\u00f76=\", \"258\u00f75=\"],\n  [\"691\u00f76=\", \"154\u00f75=\"],\n  [\"852\u00f73=\", \"425\u00f76=\"],\n  [\"239\u00f75=\", \"505\u00f78=\"],\n  [\"781\u00f78=\", \"699\u00f76=\"],\n  [\"640\u00f79=\", \"137\u00f74=\"],\n  [\"875\u00f79=\", \"499\u00f73=\"],\n  [\"141\u00f79=\", \"685\u00f73=\"],\n  [\"318\u00f78=\", \"844\u00f73=\"],\n  [\"782\u00f73=\", \"581\u00f74=\"],\n  [\"976\u00f78=\", \"978\u00f74=\"],\n  [\"850\u00f75=\", \"785\u00f72=\"],\n  [\"459\u00f76=\", \"603\u00f73=\"],\n  [\"723\u00f77=\", \"270\u00f74=\"],\n  [\"708\u00f73=\", \"492\u00f79=\"],\n  [\"722\u00f75=\", \"508\u00f76=\"],\n  [\"981\u00f75=\", \"855\u00f77=\"],\n  [\"950\u00f72=\", \"946\u00f76=\"],\n  [\"579\u00f73=\", \"420\u00f75=\"],\n  [\"343\u00f76=\", \"120\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for '\" + oldText + \"', found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NNN\u00f7N=\" division-problem prompts in the table with the\n# new set of prompts, matching each old value to its new value one-to-one.\n# Every old prompt text is unique in the document (confirmed against the\n# original OOXML), so an exact, case-sensitive whole-string Find/Replace for\n# each old value touches exactly one run.\n#\n# NOTE on ordering: \"170\u00f79=\" is both an original value (last cell) and a\n# new value (first cell), so that pair is processed first -- before any\n# other insertion creates a second \"170\u00f79=\" in the document -- so every\n# Find.Execute still targets the single, correct occurrence.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"170\u00f79=\", \"146\u00f75=\"),\n  @(\"563\u00f74=\", \"170\u00f79=\"),\n  @(\"174\u00f77=\", \"575\u00f72=\"),\n  @(\"507\u00f72=\", \"399\u00f76=\"),\n  @(\"733\u00f79=\", \"413\u00f79=\"),\n  @(\"555\u00f76=\", \"258\u00f75=\"),\n  @(\"691\u00f76=\", \"154\u00f75=\"),\n  @(\"852\u00f73=\", \"425\u00f76=\"),\n  @(\"239\u00f75=\", \"505\u00f78=\"),\n  @(\"781\u00f78=\", \"699\u00f76=\"),\n  @(\"640\u00f79=\", \"137\u00f74=\"),\n  @(\"875\u00f79=\", \"499\u00f73=\"),\n  @(\"141\u00f79=\", \"685\u00f73=\"),\n  @(\"318\u00f78=\", \"844\u00f73=\"),\n  @(\"782\u00f73=\", \"581\u00f74=\"),\n  @(\"976\u00f78=\", \"978\u00f74=\"),\n  @(\"850\u00f75=\", \"785\u00f72=\"),\n  @(\"459\u00f76=\", \"603\u00f73=\"),\n  @(\"723\u00f77=\", \"270\u00f74=\"),\n  @(\"708\u00f73=\", \"492\u00f79=\"),\n  @(\"722\u00f75=\", \"508\u00f76=\"),\n  @(\"981\u00f75=\", \"855\u00f77=\"),\n  @(\"950\u00f72=\", \"946\u00f76=\"),\n  @(\"579\u00f73=\", \"420\u00f75=\"),\n  @(\"343\u00f76=\", \"120\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"expected to find and replace '$oldText' but it was not found\"\n  }\n}\n"}
